$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.405.26"
$ws.Range("E2").Value = "  -3.78%  "
$ws.Range("D3").Value = "1.860.42"
$ws.Range("E3").Value = "  -4.46%  "
$ws.Range("E4").Value = "  -0.67%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.96"
$ws.Range("E5").Value = "  -0.17%  "
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4496"
$ws.Range("E7").Value = "  -5.58%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3863"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.67"
$ws.Range("E9").Value = "  -11.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08008"
$ws.Range("E10").Value = "  -6.27%  "
$ws.Range("E11").Value = "  -3.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.49"
$ws.Range("E12").Value = "  -2.40%  "
$ws.Range("D13").Value = "1.890.98"
$ws.Range("E13").Value = "  -2.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.900"
$ws.Range("E14").Value = "  -4.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.141"
$ws.Range("E15").Value = "  -5.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.008"
$ws.Range("E16").Value = "  -0.73%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001037"
$ws.Range("E17").Value = "  -3.72%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "86.23"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06542"
$ws.Range("E19").Value = "  -1.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.21"
$ws.Range("E20").Value = "  -7.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.007"
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.535"
$ws.Range("E22").Value = "  -4.74%  "
$ws.Range("D23").Value = "27.442.73"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.85"
$ws.Range("E24").Value = "  -5.55%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("E25").Value = "  -0.53%  "
$ws.Range("D26").Value = "2.118.25"
$ws.Range("E26").Value = "  -2.99%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "151.16"
$ws.Range("E27").Value = "  -1.75%  "
$ws.Range("E28").Value = "  -3.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.518"
$ws.Range("E29").Value = "  -6.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.034"
$ws.Range("E30").Value = "  -5.64%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.27"
$ws.Range("E31").Value = "  -1.75%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09397"
$ws.Range("E32").Value = "  -1.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.459"
$ws.Range("E33").Value = "  +1.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9292"
$ws.Range("E34").Value = "  -6.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.630"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.289"
$ws.Range("E36").Value = "  -5.32%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02231"
$ws.Range("E37").Value = "  -4.45%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.225"
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05977"
$ws.Range("E39").Value = "  -3.81%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.358"
$ws.Range("E40").Value = "  -4.25%  "
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5945"
$ws.Range("E42").Value = "  -4.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1863"
$ws.Range("E43").Value = "  -3.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.30"
$ws.Range("E44").Value = "  -6.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.278"
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5672"
$ws.Range("E46").Value = "  -4.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.23"
$ws.Range("E47").Value = "  -5.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.936"
$ws.Range("E48").Value = "  -6.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.358"
$ws.Range("E49").Value = "  -1.14%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06856"
$ws.Range("E50").Value = "  +0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.007"
